$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.131.67'
$ws.Range("E2").Value = '  -2.08%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.562.21'

$ws.Range("E4").Value = '  +0.08%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '206.31'
$ws.Range("E5").Value = '  -0.48%  '

$ws.Range("E6").Value = '  -1.89%  '

$ws.Range("E7").Value = '  +0.10%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '22.07'
$ws.Range("E8").Value = '  -0.75%  '

$ws.Range("E9").Value = '  -2.00%  '

$ws.Range("E10").Value = '  -0.04%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0861'
$ws.Range("E11").Value = '  -0.54%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.786.68'
$ws.Range("E12").Value = '  -1.60%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.569.43'
$ws.Range("E13").Value = '  -1.08%  '

$ws.Range("E14").Value = '  -2.52%  '

$ws.Range("E15").Value = '  -2.85%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '63.04'
$ws.Range("E16").Value = '  -0.74%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '27.155.84'
$ws.Range("E17").Value = '  -1.95%  '

$ws.Range("B18").Value = 'BitcoinCash'
$ws.Range("C18").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '212.27'
$ws.Range("E18").Value = '  -3.54%  '

$ws.Range("B19").Value = 'ShibaInu'
$ws.Range("C19").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0₃0688'
$ws.Range("E19").Value = '  -1.10%  '

$ws.Range("E20").Value = '  -1.73%  '

$ws.Range("E21").Value = '  +0.09%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.10'
$ws.Range("E22").Value = '  -1.09%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.39'
$ws.Range("E23").Value = '  -2.14%  '

$ws.Range("E24").Value = '  +0.38%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '152.29'
$ws.Range("E25").Value = '  -0.58%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.61'
$ws.Range("E26").Value = '  -3.65%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '14.84'
$ws.Range("E27").Value = '  -2.05%  '

$ws.Range("E29").Value = '  -2.04%  '

$ws.Range("E30").Value = '  -1.04%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0463'
$ws.Range("E31").Value = '  -1.21%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.16'
$ws.Range("E32").Value = '  -1.86%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.375.37'
$ws.Range("E33").Value = '  +0.29%  '

$ws.Range("E34").Value = '  +0.65%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.54'
$ws.Range("E35").Value = '  +0.58%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.945'
$ws.Range("E37").Value = '  -3.32%  '

$ws.Range("E38").Value = '  -1.58%  '

$ws.Range("E39").Value = '  -3.18%  '

$ws.Range("E40").Value = '  -1.39%  '

$ws.Range("E41").Value = '  +0.06%  '

$ws.Range("E42").Value = '  +1.45%  '

$ws.Range("E43").Value = '  +3.96%  '

$ws.Range("E44").Value = '  +0.00%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '63.39'
$ws.Range("E45").Value = '  -1.40%  '

$ws.Range("E46").Value = '  -0.88%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.698.90'
$ws.Range("E47").Value = '  -1.58%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '85.45'
$ws.Range("E48").Value = '  -2.64%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0⁷0995'
$ws.Range("E49").Value = '  -1.20%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0491'
$ws.Range("E50").Value = '  -1.08%  '

$ws.Range("E51").Value = '  +0.11%  '
